$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the row-level edits from the updated export, processed from the
# bottom of the sheet upward so that earlier (lower) row numbers used below
# stay valid while later (higher-numbered) operations are applied first.

# Insert 1 new row(s) after row 163
$ws.Rows.Item(164).Insert()
$ws.Cells.Item(164, 1).Value = "'004208447"
$ws.Cells.Item(164, 2).Value = "LEILA"
$ws.Cells.Item(164, 3).Value = 39.86

# Remove row 150 (balance no longer present in the new export)
$ws.Rows.Item(150).Delete()

# Remove row 64 (balance no longer present in the new export)
$ws.Rows.Item(64).Delete()

# Remove row 57 (balance no longer present in the new export)
$ws.Rows.Item(57).Delete()

# Remove row 46 (balance no longer present in the new export)
$ws.Rows.Item(46).Delete()

# Insert 3 new row(s) after row 38
$ws.Rows.Item(39).Resize(3).Insert()
$ws.Cells.Item(39, 1).Value = "'004260002"
$ws.Cells.Item(39, 2).Value = "ERICA"
$ws.Cells.Item(39, 3).Value = 346.75
$ws.Cells.Item(40, 1).Value = "'004240400"
$ws.Cells.Item(40, 2).Value = "ADRIANA"
$ws.Cells.Item(40, 3).Value = 299.85
$ws.Cells.Item(41, 1).Value = "'001651617"
$ws.Cells.Item(41, 2).Value = "MIRELLA"
$ws.Cells.Item(41, 3).Value = 291.75

# Insert 1 new row(s) after row 26
$ws.Rows.Item(27).Insert()
$ws.Cells.Item(27, 1).Value = "'003553997"
$ws.Cells.Item(27, 2).Value = "MIRELLA"
$ws.Cells.Item(27, 3).Value = 776.2

# Insert 2 new row(s) after row 19
$ws.Rows.Item(20).Resize(2).Insert()
$ws.Cells.Item(20, 1).Value = "'004261201"
$ws.Cells.Item(20, 2).Value = "ANA"
$ws.Cells.Item(20, 3).Value = 1182.02
$ws.Cells.Item(21, 1).Value = "'004210959"
$ws.Cells.Item(21, 2).Value = "ANA"
$ws.Cells.Item(21, 3).Value = 1176.16

# Insert 1 new row(s) after row 18
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = "'004322719"
$ws.Cells.Item(19, 2).Value = "GISELA"
$ws.Cells.Item(19, 3).Value = 1815.24

# Remove row 17 (balance no longer present in the new export)
$ws.Rows.Item(17).Delete()

# Update row 14 in place (balance refreshed, account changed)
$ws.Cells.Item(14, 1).Value = "'004363260"
$ws.Cells.Item(14, 2).Value = "LARISSA"
$ws.Cells.Item(14, 3).Value = 10045.41

Write-Output "edits applied"